# Commit: feat: add 2022-Q3 data
#
# This adds a new "2022-Q3" worksheet (placed right after the "总计"
# summary sheet) with the quarter's fund-holding detail table, and
# records its headline numbers (9 holdings, 0.63 billion yuan) as a new
# first data row on the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计" (i.e. before the
#    sheet that is currently "2022-Q1", which sits in tab position 2).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$oldQ1   = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($oldQ1, $null)
$newSheet.Name = "2022-Q3"

# Re-fetch sheet references now that the collection has shifted - a
# reference obtained before the Add() call can paste/copy against a
# stale position.
$summary = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Item(2)
$oldQ1 = $wb.Worksheets.Item(3)

# Match the page margins used throughout the rest of the workbook
# (0.75in / 0.75in / 1in / 1in / 0.5in / 0.5in -> points).
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 2. Fill in the header row (B1:H1) and copy the header style from the
#    equivalent cells on the neighbouring "2022-Q1" sheet so the new
#    sheet matches the workbook's existing look (bold, centred, boxed).
# ---------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$oldQ1.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Fill in the fund detail rows 2-10.
#    Column A: numeric row index (0-based)
#    Column B: fund code (text - preserve any leading zeros)
#    Column C: fund name (text)
#    Columns D-G: text-formatted numbers (kept as text, matching the
#                 rest of the workbook's quarterly tables)
#    Column H: numeric ranking
# ---------------------------------------------------------------------
function Set-TextCell($cell, $text) {
    $cell.Value = "'" + $text
}

$rows = @(
    @{ A=0; B="009562"; C="工银全球股票（QDII）美元";                        D="5.89"; E="93.72"; F="2.12"; G="0.1249"; H=7 },
    @{ A=1; B="009563"; C="工银全球股票（QDII）港币";                        D="5.89"; E="93.72"; F="2.12"; G="0.1249"; H=7 },
    @{ A=2; B="486001"; C="工银瑞信中国机会全球配置股票（QDII）人民币";      D="5.89"; E="93.72"; F="2.12"; G="0.1249"; H=7 },
    @{ A=3; B="486002"; C="工银全球精选股票（QDII）";                        D="3.72"; E="93.69"; F="2.96"; G="0.1101"; H=3 },
    @{ A=4; B="009225"; C="天弘中证中美互联网指数（QDII）A";                 D="1.20"; E="94.98"; F="4.47"; G="0.0536"; H=9 },
    @{ A=5; B="009226"; C="天弘中证中美互联网指数（QDII）C";                 D="0.60"; E="94.98"; F="4.47"; G="0.0268"; H=9 },
    @{ A=6; B="012751"; C="建信纳斯达克100指数（QDII）A 美元现汇";           D="0.64"; E="80.13"; F="3.48"; G="0.0223"; H=5 },
    @{ A=7; B="012752"; C="建信纳斯达克100指数（QDII）C 人民币";             D="0.64"; E="80.13"; F="3.48"; G="0.0223"; H=5 },
    @{ A=8; B="012753"; C="建信纳斯达克100指数（QDII）C 美元现汇";           D="0.64"; E="80.13"; F="3.48"; G="0.0223"; H=5 }
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value2 = $row.A
    Set-TextCell $newSheet.Range("B$r") $row.B
    $newSheet.Range("C$r").Value = $row.C
    Set-TextCell $newSheet.Range("D$r") $row.D
    Set-TextCell $newSheet.Range("E$r") $row.E
    Set-TextCell $newSheet.Range("F$r") $row.F
    Set-TextCell $newSheet.Range("G$r") $row.G
    $newSheet.Range("H$r").Value2 = $row.H
    $r = $r + 1
}

# Copy the "A" index-column style from the neighbouring sheet too.
$oldQ1.Range("A2").Copy()
$newSheet.Range("A2:A10").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4. Update the "总计" summary sheet: push the existing three rows down
#    by one and insert the new 2022-Q3 totals as the new row 2. Values
#    are written directly (not via Rows.Insert) to avoid floating point
#    re-encoding drift.
# ---------------------------------------------------------------------
$summary.Range("A5").Value2 = 3
$summary.Range("B5").Value  = "2021-Q2"
$summary.Range("C5").Value2 = 7
$summary.Range("D5").Value2 = 1.04

$summary.Range("A4").Value2 = 2
$summary.Range("B4").Value  = "2021-Q4"
$summary.Range("C4").Value2 = 9
$summary.Range("D4").Value2 = 2.33

$summary.Range("A3").Value2 = 1
$summary.Range("B3").Value  = "2022-Q1"
$summary.Range("C3").Value2 = 9
$summary.Range("D3").Value2 = 0.92

$summary.Range("A2").Value2 = 0
$summary.Range("B2").Value  = "2022-Q3"
$summary.Range("C2").Value2 = 9
$summary.Range("D2").Value2 = 0.63

# Make sure the new A5 row carries the same index-column style as the
# rows above it.
$summary.Range("A4").Copy()
$summary.Range("A5").PasteSpecial(-4122)

# Restore "总计" as the active sheet (matches the workbook's unchanged
# bookViews/activeTab="0").
$summary.Activate()
